$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

$ws.Range("D5").Value = "pass"
$ws.Range("E5").Value = "Videos"
$ws.Range("F5").Value = "https://timesofindia.indiatimes.com/entertainment/bengali/movie-reviews/macher-jhol/movie-review/60174378.cms"
